$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.526.17'
$ws.Range("E2").Value = '  +3.92%  '

$ws.Range("D3").Value = '2.432.24'
$ws.Range("E3").Value = '  +2.82%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.77'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.76'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.52%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.515'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.20%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("E9").Value = '  +5.51%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.43'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.40%  '

$ws.Range("E11").Value = '  +1.89%  '

$ws.Range("E12").Value = '  +1.12%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.85'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.97%  '

$ws.Range("E14").Value = '  +3.44%  '

$ws.Range("D15").Value = '2.811.57'
$ws.Range("E15").Value = '  +2.77%  '

$ws.Range("D16").Value = '2.461.17'
$ws.Range("E16").Value = '  +3.91%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.840'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.19%  '

$ws.Range("D18").Value = '44.483.22'
$ws.Range("E18").Value = '  +3.92%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.39'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.70%  '

$ws.Range("E20").Value = '  +2.18%  '

$ws.Range("D21").Value = '0.0₃0908'
$ws.Range("E21").Value = '  +2.58%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.92'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.52%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '241.46'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.81%  '

$ws.Range("E24").Value = '  +4.48%  '

$ws.Range("E25").Value = '  +2.38%  '

$ws.Range("E26").Value = '  -0.03%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.23'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.09%  '

$ws.Range("E28").Value = '  -4.18%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.66'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.80%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.44'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.41%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '48.61'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.53%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.123'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +17.42%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.49'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +12.31%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.21'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.61%  '

$ws.Range("E35").Value = '  +0.25%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0767'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.07%  '

$ws.Range("E37").Value = '  +2.81%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.54'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.41%  '

$ws.Range("E39").Value = '  +4.89%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '125.67'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.39%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.22'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.77%  '

$ws.Range("E42").Value = '  +1.26%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.62'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.24%  '

$ws.Range("E44").Value = '  +3.90%  '

$ws.Range("D45").Value = '1.947.93'
$ws.Range("E45").Value = '  +0.81%  '

$ws.Range("E46").Value = '  +2.32%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.96'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +8.72%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.80'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.96%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.69'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +11.53%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.58'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.50%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.91'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.54%  '

